# Sprint2Tema1Ejercicios.xlsx - "Up to excercise 10"
#
# The two exercises about listing manufacturer codes that have products
# ("Lista el codigo de los fabricantes que tienen productos en la ... ")
# were re-worded (mesa -> tabla) and re-positioned to rows 12-13 (right
# after exercise 10/11 about truncating/rounding prices), pushing the
# "ascendente" exercise (and everything after it) down to keep a clean
# sequential numbering. The worksheet's row numbers in column B are driven
# by formulas and do not change; only the description text in column C
# for rows 12 and 13 actually needs new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ejercicios")

$ws.Range("C12").Value = "Lista el código de los fabricantes que tienen productos en la tabla ""producto""."
$ws.Range("C13").Value = "Lista el código de los fabricantes que tienen productos en la tabla ""producto"", eliminando los códigos que aparecen repetidos."

# Move the active cell / selection to C13, matching the saved cursor
# position recorded in the workbook.
[void]$ws.Range("C13").Select()
